# Updates the cryptocurrency price (D) and 1h volume change (E) columns
# for rows 2-51 on Sheet1, matching the refreshed GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "20.106.51"
$ws.Range("E2").Value = "  -1.49%  "
$ws.Range("D3").Value = "1.423.10"
$ws.Range("E3").Value = "  -1.39%  "
$ws.Range("D4").NumberFormat = "@"  # keep numeric-looking price as literal text, like the source data
$ws.Range("D4").Value = "0.9997"
$ws.Range("E4").Value = "  -0.51%  "
$ws.Range("D5").NumberFormat = "@"  # keep numeric-looking price as literal text, like the source data
$ws.Range("D5").Value = "0.9971"
$ws.Range("E5").Value = "  -0.70%  "
$ws.Range("D6").NumberFormat = "@"  # keep numeric-looking price as literal text, like the source data
$ws.Range("D6").Value = "277.03"
$ws.Range("E6").Value = "  -0.29%  "
$ws.Range("D7").NumberFormat = "@"  # keep numeric-looking price as literal text, like the source data
$ws.Range("D7").Value = "0.3709"
$ws.Range("E7").Value = "  -0.96%  "
$ws.Range("D8").NumberFormat = "@"  # keep numeric-looking price as literal text, like the source data
$ws.Range("D8").Value = "0.3153"
$ws.Range("E8").Value = "  +2.67%  "
$ws.Range("D9").NumberFormat = "@"  # keep numeric-looking price as literal text, like the source data
$ws.Range("D9").Value = "39.45"
$ws.Range("E9").Value = "  -2.75%  "
$ws.Range("D10").NumberFormat = "@"  # keep numeric-looking price as literal text, like the source data
$ws.Range("D10").Value = "1.060"
$ws.Range("E10").Value = "  +4.62%  "
$ws.Range("D11").NumberFormat = "@"  # keep numeric-looking price as literal text, like the source data
$ws.Range("D11").Value = "0.06550"
$ws.Range("E11").Value = "  +0.04%  "
$ws.Range("D12").NumberFormat = "@"  # keep numeric-looking price as literal text, like the source data
$ws.Range("D12").Value = "0.9973"
$ws.Range("E12").Value = "  -0.75%  "
$ws.Range("D13").NumberFormat = "@"  # keep numeric-looking price as literal text, like the source data
$ws.Range("D13").Value = "5.530"
$ws.Range("E13").Value = "  +2.73%  "
$ws.Range("D14").NumberFormat = "@"  # keep numeric-looking price as literal text, like the source data
$ws.Range("D14").Value = "17.95"
$ws.Range("E14").Value = "  +3.96%  "
$ws.Range("D15").NumberFormat = "@"  # keep numeric-looking price as literal text, like the source data
$ws.Range("D15").Value = "6.205"
$ws.Range("E15").Value = "  +1.16%  "
$ws.Range("D16").Value = "1.421.90"
$ws.Range("E16").Value = "  -1.57%  "
$ws.Range("D17").NumberFormat = "@"  # keep numeric-looking price as literal text, like the source data
$ws.Range("D17").Value = "0.00001024"
$ws.Range("E17").Value = "  +1.49%  "
$ws.Range("D18").NumberFormat = "@"  # keep numeric-looking price as literal text, like the source data
$ws.Range("D18").Value = "0.05714"
$ws.Range("E18").Value = "  -2.77%  "
$ws.Range("D19").NumberFormat = "@"  # keep numeric-looking price as literal text, like the source data
$ws.Range("D19").Value = "0.9973"
$ws.Range("E19").Value = "  -0.68%  "
$ws.Range("D20").NumberFormat = "@"  # keep numeric-looking price as literal text, like the source data
$ws.Range("D20").Value = "71.84"
$ws.Range("E20").Value = "  -5.65%  "
$ws.Range("D21").NumberFormat = "@"  # keep numeric-looking price as literal text, like the source data
$ws.Range("D21").Value = "5.622"
$ws.Range("E21").Value = "  -1.99%  "
$ws.Range("D22").NumberFormat = "@"  # keep numeric-looking price as literal text, like the source data
$ws.Range("D22").Value = "14.89"
$ws.Range("E22").Value = "  +3.62%  "
$ws.Range("D23").NumberFormat = "@"  # keep numeric-looking price as literal text, like the source data
$ws.Range("D23").Value = "11.11"
$ws.Range("E23").Value = "  +1.98%  "
$ws.Range("D24").NumberFormat = "@"  # keep numeric-looking price as literal text, like the source data
$ws.Range("D24").Value = "2.227"
$ws.Range("E24").Value = "  -3.44%  "
$ws.Range("D25").Value = "20.131.28"
$ws.Range("E25").Value = "  -1.38%  "
$ws.Range("D26").NumberFormat = "@"  # keep numeric-looking price as literal text, like the source data
$ws.Range("D26").Value = "2.300"
$ws.Range("E26").Value = "  +3.75%  "
$ws.Range("D27").NumberFormat = "@"  # keep numeric-looking price as literal text, like the source data
$ws.Range("D27").Value = "134.49"
$ws.Range("E27").Value = "  -6.39%  "
$ws.Range("D28").NumberFormat = "@"  # keep numeric-looking price as literal text, like the source data
$ws.Range("D28").Value = "17.32"
$ws.Range("E28").Value = "  +1.75%  "
$ws.Range("D29").Value = "1.581.15"
$ws.Range("E29").Value = "  -1.63%  "
$ws.Range("D30").NumberFormat = "@"  # keep numeric-looking price as literal text, like the source data
$ws.Range("D30").Value = "111.05"
$ws.Range("E30").Value = "  +1.59%  "
$ws.Range("D31").NumberFormat = "@"  # keep numeric-looking price as literal text, like the source data
$ws.Range("D31").Value = "3.960"
$ws.Range("E31").Value = "  +7.55%  "
$ws.Range("D32").NumberFormat = "@"  # keep numeric-looking price as literal text, like the source data
$ws.Range("D32").Value = "5.297"
$ws.Range("E32").Value = "  -2.29%  "
$ws.Range("D33").NumberFormat = "@"  # keep numeric-looking price as literal text, like the source data
$ws.Range("D33").Value = "0.8281"
$ws.Range("E33").Value = "  -8.17%  "
$ws.Range("D34").NumberFormat = "@"  # keep numeric-looking price as literal text, like the source data
$ws.Range("D34").Value = "0.07808"
$ws.Range("E34").Value = "  +0.97%  "
$ws.Range("D35").NumberFormat = "@"  # keep numeric-looking price as literal text, like the source data
$ws.Range("D35").Value = "1.477"
$ws.Range("E35").Value = "  +6.39%  "
$ws.Range("D36").NumberFormat = "@"  # keep numeric-looking price as literal text, like the source data
$ws.Range("D36").Value = "4.932"
$ws.Range("E36").Value = "  +4.22%  "
$ws.Range("D37").NumberFormat = "@"  # keep numeric-looking price as literal text, like the source data
$ws.Range("D37").Value = "0.05870"
$ws.Range("E37").Value = "  +4.45%  "
$ws.Range("D38").NumberFormat = "@"  # keep numeric-looking price as literal text, like the source data
$ws.Range("D38").Value = "8.027"
$ws.Range("E38").Value = "  -3.05%  "
$ws.Range("D39").NumberFormat = "@"  # keep numeric-looking price as literal text, like the source data
$ws.Range("D39").Value = "0.9965"
$ws.Range("E39").Value = "  -0.62%  "
$ws.Range("D40").NumberFormat = "@"  # keep numeric-looking price as literal text, like the source data
$ws.Range("D40").Value = "10.63"
$ws.Range("E40").Value = "  -2.09%  "
$ws.Range("D41").NumberFormat = "@"  # keep numeric-looking price as literal text, like the source data
$ws.Range("D41").Value = "0.02069"
$ws.Range("E41").Value = "  +1.53%  "
$ws.Range("D42").NumberFormat = "@"  # keep numeric-looking price as literal text, like the source data
$ws.Range("D42").Value = "1.112"
$ws.Range("E42").Value = "  -2.83%  "
$ws.Range("D43").NumberFormat = "@"  # keep numeric-looking price as literal text, like the source data
$ws.Range("D43").Value = "0.1877"
$ws.Range("E43").Value = "  -1.96%  "
$ws.Range("D44").NumberFormat = "@"  # keep numeric-looking price as literal text, like the source data
$ws.Range("D44").Value = "0.5350"
$ws.Range("E44").Value = "  +0.67%  "
$ws.Range("D45").NumberFormat = "@"  # keep numeric-looking price as literal text, like the source data
$ws.Range("D45").Value = "12.36"
$ws.Range("E45").Value = "  +2.31%  "
$ws.Range("D46").NumberFormat = "@"  # keep numeric-looking price as literal text, like the source data
$ws.Range("D46").Value = "3.553"
$ws.Range("E46").Value = "  -1.03%  "
$ws.Range("D47").NumberFormat = "@"  # keep numeric-looking price as literal text, like the source data
$ws.Range("D47").Value = "118.32"
$ws.Range("E47").Value = "  +6.03%  "
$ws.Range("D48").NumberFormat = "@"  # keep numeric-looking price as literal text, like the source data
$ws.Range("D48").Value = "0.5242"
$ws.Range("E48").Value = "  +1.67%  "
$ws.Range("D49").NumberFormat = "@"  # keep numeric-looking price as literal text, like the source data
$ws.Range("D49").Value = "1.785"
$ws.Range("E49").Value = "  -0.21%  "
$ws.Range("D50").NumberFormat = "@"  # keep numeric-looking price as literal text, like the source data
$ws.Range("D50").Value = "1.041"
$ws.Range("E50").Value = "  -1.15%  "
$ws.Range("D51").NumberFormat = "@"  # keep numeric-looking price as literal text, like the source data
$ws.Range("D51").Value = "0.9971"
$ws.Range("E51").Value = "  -0.58%  "
